# README_CheatSheet content is now integrated into the main README: the
# workbook gains a second "Analysis" worksheet (the old summary metrics +
# new URL/Dimension failure breakdown columns), while the first worksheet
# ("Test Results") is rewritten to hold the per-test-case results table.

$wb = $excel.ActiveWorkbook

# --- existing sheet becomes "Test Results" test-case table ---------------
$wsResults = $wb.Worksheets.Item(1)
$wsResults.Name = "Test Results"

# Clear out the old Metric/Value summary content first.
$wsResults.Cells.Clear()

$resultsHeaders = @("URL", "Parameter", "Result", "Details")
for ($col = 1; $col -le $resultsHeaders.Length; $col++) {
    $cell = $wsResults.Cells.Item(1, $col)
    $cell.Value = $resultsHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.Font.Color = 16777215
    $cell.Interior.Color = 12419407
}

$resultsRows = @(
    @("https://example.com/test1", "param1", "Pass", "Details of test case 1"),
    @("https://example.com/test2", "param2", "Fail", "Details of test case 2"),
    @("https://example.com/test3", "param3", "Pass", "Details of test case 3"),
    @("https://example.com/test4", "param4", "Fail", "Details of test case 4")
)

for ($r = 0; $r -lt $resultsRows.Length; $r++) {
    $row = $resultsRows[$r]
    for ($col = 1; $col -le $row.Length; $col++) {
        $wsResults.Cells.Item($r + 2, $col).Value = $row[$col - 1]
    }
}

# --- new "Analysis" sheet, placed right after "Test Results" -------------
$wsAnalysis = $wb.Worksheets.Add($null, $wsResults)
$wsAnalysis.Name = "Analysis"

# Columns 3 and 6 are blank spacer cells between the three mini-tables,
# but they still carry the header formatting (mirrors the source row).
$analysisHeaders = @{
    1 = "Metric";           2 = "Value";  3 = $null
    4 = "URL Failures";     5 = "Count";  6 = $null
    7 = "Dimension Failures"; 8 = "Count"
}
for ($col = 1; $col -le 8; $col++) {
    $cell = $wsAnalysis.Cells.Item(1, $col)
    $headerText = $analysisHeaders[$col]
    if ($headerText) {
        $cell.Value = $headerText
    }
    $cell.Font.Bold = $true
    $cell.Font.Color = 16777215
    $cell.Interior.Color = 12419407
}

$metricRows = @(
    @("Total Calls", 1),
    @("Total Test Cases", 1),
    @("Total Passes", 1),
    @("Total Fails", 1),
    @("Pass Percentage", "100%"),
    @("Fail Percentage", "0%")
)
for ($r = 0; $r -lt $metricRows.Length; $r++) {
    $row = $metricRows[$r]
    $valueCell = $wsAnalysis.Cells.Item($r + 2, 2)
    $wsAnalysis.Cells.Item($r + 2, 1).Value = $row[0]
    # Keep strings like "100%"/"0%" as literal text instead of letting
    # COM auto-coerce them into a percentage number + format.
    if ($row[1] -is [string]) {
        $valueCell.NumberFormat = "@"
    }
    $valueCell.Value = $row[1]
}

$urlFailureRows = @(
    @("https://example.com/failure1", 5),
    @("https://example.com/failure2", 3)
)
for ($r = 0; $r -lt $urlFailureRows.Length; $r++) {
    $row = $urlFailureRows[$r]
    $wsAnalysis.Cells.Item($r + 2, 4).Value = $row[0]
    $wsAnalysis.Cells.Item($r + 2, 5).Value = $row[1]
}

$dimensionFailureRows = @(
    @("example_dimension_1", 2),
    @("example_dimension_2", 4)
)
for ($r = 0; $r -lt $dimensionFailureRows.Length; $r++) {
    $row = $dimensionFailureRows[$r]
    $wsAnalysis.Cells.Item($r + 2, 7).Value = $row[0]
    $wsAnalysis.Cells.Item($r + 2, 8).Value = $row[1]
}

$wsResults.Activate()
